# Add two new rows (Anapa & Gelendzhik resort towns) as new local extremums
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row19 = @(3703000, "город-курорт Анапа", 2021, 1.749007936507937, 1.421437873576652, 0.26429432319309099, 0.48113729915185383, 0.35689377233967179, 1.3045617637456459, 0.85349502337851924, 0.62641631595985758, 0.19801656001269319, 0.1659048139950032, [double]"5.6246756347337092E-2", [double]"4.6299740287124642E-2", [double]"8.4231657304698185E-2", [double]"3.9654239413479959E-2", [double]"5.1819450382683742E-2", 0.31759183029793209, 0.2386452468783129)

$row20 = @(3708000, "город-курорт Геленджик", 2021, [double]"-4.265873015873016E-2", 0.76834407453167264, 0.32300945528313751, 0.52922955069295663, 0.68402772370909792, 2.82202011720989, 1.0403368716792889, 0.69116218841048882, 0.23877624403126041, 0.66356146232570201, [double]"4.296398496954567E-2", [double]"6.1342243182431729E-3", [double]"1.2013378495212151E-2", [double]"8.8894538935417808E-3", [double]"5.7063138863778722E-2", 0.3936003802977765, 0.13065006111681071)

for ($i = 0; $i -lt $row19.Length; $i++) {
    $cell = $ws.Cells.Item(19, $i + 1)
    $cell.Value = $row19[$i]
    $cell.HorizontalAlignment = -4108
}

for ($i = 0; $i -lt $row20.Length; $i++) {
    $cell = $ws.Cells.Item(20, $i + 1)
    $cell.Value = $row20[$i]
    $cell.HorizontalAlignment = -4108
}
